$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A59").Value = "Metodología mides"
$ws.Range("B59").Value = "metodologia_mides"
